{"js": "// Apply proofing-mark / spell-check markup edits described by the diff:\n//  - Title run split into LoRa / PodWave proofErr-wrapped runs\n//  - GetPODMeta / GetPODData / GetPODResult paragraphs get spellStart/\n//    gramStart/spellEnd/gramEnd proofErr wrapping and a literal \"()\" suffix\n//  - SendPODData paragraph gets spellStart/spellEnd wrapping and a literal\n//    \"(String)\" suffix\n//  - UpdatePODParameter paragraph gets spellStart/spellEnd wrapping and the\n//    \"_GoBack\" bookmark (previously living in its own empty paragraph a few\n//    paragraphs down) moved to sit right after its run\n//  - the paragraph that used to hold the \"_GoBack\" bookmark becomes empty\n\nconst OOXML_NS =\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">';\n\nfunction pkg(bodyInner) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    OOXML_NS +\n    \"<w:body>\" +\n    bodyInner +\n    \"</w:body></w:document>\" +\n    \"</pkg:xmlData></pkg:part></pkg:package>\"\n  );\n}\n\n// Builds replacement XML for a \"GetPODMeta()\"-style paragraph: name wrapped\n// in spellStart/gramStart .. spellEnd, then a literal \"()\" run, then gramEnd.\nfunction callParagraphXml(styleVal, name) {\n  return pkg(\n    \"<w:p>\" +\n      '<w:pPr><w:pStyle w:val=\"' +\n      styleVal +\n      '\"/></w:pPr>' +\n      '<w:proofErr w:type=\"spellStart\"/>' +\n      '<w:proofErr w:type=\"gramStart\"/>' +\n      \"<w:r><w:t>\" +\n      name +\n      \"</w:t></w:r>\" +\n      '<w:proofErr w:type=\"spellEnd\"/>' +\n      \"<w:r><w:t>()</w:t></w:r>\" +\n      '<w:proofErr w:type=\"gramEnd\"/>' +\n      \"</w:p>\"\n  );\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text,style\");\nawait context.sync();\n\nfunction findByText(text) {\n  for (const p of paragraphs.items) {\n    if (p.text === text) return p;\n  }\n  throw new Error('Paragraph with text \"' + text + '\" not found');\n}\n\n// 1) Title paragraph \u2014 split \"LoRa PodWave Wifi-Server API Documentation :\"\n// NB: the space right before the trailing colon is a non-breaking space\n// (U+00A0) in the source document, matching the same French-typography\n// convention used by the other headings in this doc (\"Pr\u00e9ambule ... :\",\n// \"API Routes : \", etc.) \u2014 keep it byte-for-byte identical.\nconst NBSP = \"\\u00A0\";\nconst titlePara = findByText(\n  \"LoRa PodWave Wifi-Server API Documentation\" + NBSP + \":\"\n);\nconst titleXml = pkg(\n  \"<w:p>\" +\n    '<w:pPr><w:pStyle w:val=\"Title\"/></w:pPr>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    \"<w:r><w:t>LoRa</w:t></w:r>\" +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    \"<w:r><w:t>PodWave</w:t></w:r>\" +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> Wifi-Server API Documentation' +\n    NBSP +\n    ':</w:t></w:r>' +\n    \"</w:p>\"\n);\ntitlePara.insertOoxml(titleXml, Word.InsertLocation.replace);\n\n// 2) GetPODMeta / GetPODData / GetPODResult -> add \"()\" + spell/gram proofErr\nfindByText(\"GetPODMeta\").insertOoxml(\n  callParagraphXml(\"Heading3\", \"GetPODMeta\"),\n  Word.InsertLocation.replace\n);\nfindByText(\"GetPODData\").insertOoxml(\n  callParagraphXml(\"Heading3\", \"GetPODData\"),\n  Word.InsertLocation.replace\n);\nfindByText(\"GetPODResult\").insertOoxml(\n  callParagraphXml(\"Heading3\", \"GetPODResult\"),\n  Word.InsertLocation.replace\n);\n\n// 3) SendPODData -> add \"(String)\" + spellStart/spellEnd proofErr only\nconst sendXml = pkg(\n  \"<w:p>\" +\n    '<w:pPr><w:pStyle w:val=\"Heading3\"/></w:pPr>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    \"<w:r><w:t>SendPODData</w:t></w:r>\" +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    \"<w:r><w:t>(String)</w:t></w:r>\" +\n    \"</w:p>\"\n);\nfindByText(\"SendPODData\").insertOoxml(sendXml, Word.InsertLocation.replace);\n\n// 4) UpdatePODParameter -> wrap in spellStart/spellEnd and pull the\n//    \"_GoBack\" bookmark in right after the run\nconst updateXml = pkg(\n  \"<w:p>\" +\n    '<w:pPr><w:pStyle w:val=\"Heading3\"/></w:pPr>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    \"<w:r><w:t>UpdatePODParameter</w:t></w:r>\" +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n    '<w:bookmarkEnd w:id=\"0\"/>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    \"</w:p>\"\n);\nfindByText(\"UpdatePODParameter\").insertOoxml(\n  updateXml,\n  Word.InsertLocation.replace\n);\n\nawait context.sync();\n\n// 5) The paragraph a few rows below UpdatePODParameter that used to hold\n//    the \"_GoBack\" bookmark becomes a plain empty paragraph. Locate it by\n//    walking forward from UpdatePODParameter (it's the 3rd empty paragraph\n//    after it) and clear it with an empty-paragraph OOXML replace so any\n//    leftover bookmark markup is gone.\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nlet updateIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"UpdatePODParameter\") {\n    updateIndex = i;\n    break;\n  }\n}\nif (updateIndex === -1) {\n  throw new Error(\"UpdatePODParameter paragraph not found on second pass\");\n}\n\nlet emptySeen = 0;\nlet bookmarkParaIndex = -1;\nfor (let i = updateIndex + 1; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"\") {\n    emptySeen++;\n    if (emptySeen === 3) {\n      bookmarkParaIndex = i;\n      break;\n    }\n  } else {\n    break;\n  }\n}\n\nif (bookmarkParaIndex !== -1) {\n  const emptyXml = pkg(\"<w:p/>\");\n  paragraphs.items[bookmarkParaIndex].insertOoxml(\n    emptyXml,\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n", "ps1": "# Apply proofing-mark / spell-check markup edits described by the diff:\n#  - Title run split into LoRa / PodWave proofErr-wrapped runs\n#  - GetPODMeta / GetPODData / GetPODResult paragraphs get spellStart/\n#    gramStart/spellEnd/gramEnd proofErr wrapping and a literal \"()\" suffix\n#  - SendPODData paragraph gets spellStart/spellEnd wrapping and a literal\n#    \"(String)\" suffix\n#  - UpdatePODParameter paragraph gets spellStart/spellEnd wrapping and the\n#    \"_GoBack\" bookmark (previously living in its own empty paragraph a few\n#    paragraphs down) moved to sit right after its run\n#  - the paragraph that used to hold the \"_GoBack\" bookmark becomes empty\n\n$d = $word.ActiveDocument\n\n$xmlHeader = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>'\n$xmlFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\nfunction New-Pkg([string]$bodyInner) {\n    return $xmlHeader + $bodyInner + $xmlFooter\n}\n\nfunction New-CallParagraphXml([string]$styleVal, [string]$name) {\n    $inner = '<w:p><w:pPr><w:pStyle w:val=\"' + $styleVal + '\"/></w:pPr>' +\n        '<w:proofErr w:type=\"spellStart\"/><w:proofErr w:type=\"gramStart\"/>' +\n        '<w:r><w:t>' + $name + '</w:t></w:r>' +\n        '<w:proofErr w:type=\"spellEnd\"/>' +\n        '<w:r><w:t>()</w:t></w:r>' +\n        '<w:proofErr w:type=\"gramEnd\"/></w:p>'\n    return New-Pkg $inner\n}\n\nfunction Find-ParagraphByText([string]$text) {\n    foreach ($p in $d.Paragraphs) {\n        $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n        if ($t -eq $text) {\n            return $p\n        }\n    }\n    return $null\n}\n\n# 1) Title paragraph - split \"LoRa PodWave Wifi-Server API Documentation :\"\n# NB: the space right before the trailing colon is a non-breaking space\n# (U+00A0) in the source document, matching the same French-typography\n# convention used by the other headings in this doc (\"Pr\u00e9ambule ... :\",\n# \"API Routes : \", etc.) - keep it byte-for-byte identical.\n$nbsp = [char]0x00A0\n$titleText = \"LoRa PodWave Wifi-Server API Documentation\" + $nbsp + \":\"\n$titlePara = Find-ParagraphByText $titleText\nif ($null -eq $titlePara) { throw \"Title paragraph not found\" }\n$titleInner = '<w:p><w:pPr><w:pStyle w:val=\"Title\"/></w:pPr>' +\n    '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>LoRa</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>PodWave</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> Wifi-Server API Documentation' + $nbsp + ':</w:t></w:r>' +\n    '</w:p>'\n[void]$titlePara.Range.InsertXML((New-Pkg $titleInner))\n\n# 2) GetPODMeta / GetPODData / GetPODResult -> add \"()\" + spell/gram proofErr\nforeach ($name in @(\"GetPODMeta\", \"GetPODData\", \"GetPODResult\")) {\n    $para = Find-ParagraphByText $name\n    if ($null -eq $para) { throw \"$name paragraph not found\" }\n    [void]$para.Range.InsertXML((New-CallParagraphXml \"Heading3\" $name))\n}\n\n# 3) SendPODData -> add \"(String)\" + spellStart/spellEnd proofErr only\n$sendPara = Find-ParagraphByText \"SendPODData\"\nif ($null -eq $sendPara) { throw \"SendPODData paragraph not found\" }\n$sendInner = '<w:p><w:pPr><w:pStyle w:val=\"Heading3\"/></w:pPr>' +\n    '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>SendPODData</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t>(String)</w:t></w:r>' +\n    '</w:p>'\n[void]$sendPara.Range.InsertXML((New-Pkg $sendInner))\n\n# 4) UpdatePODParameter -> wrap in spellStart/spellEnd and pull the\n#    \"_GoBack\" bookmark in right after the run\n$updatePara = Find-ParagraphByText \"UpdatePODParameter\"\nif ($null -eq $updatePara) { throw \"UpdatePODParameter paragraph not found\" }\n$updateInner = '<w:p><w:pPr><w:pStyle w:val=\"Heading3\"/></w:pPr>' +\n    '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>UpdatePODParameter</w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '</w:p>'\n[void]$updatePara.Range.InsertXML((New-Pkg $updateInner))\n\n# 5) The paragraph a few rows below UpdatePODParameter that used to hold the\n#    \"_GoBack\" bookmark becomes a plain empty paragraph. Locate it by walking\n#    forward from UpdatePODParameter (it's the 3rd empty paragraph after it)\n#    and clear it with an empty-paragraph OOXML replace so any leftover\n#    bookmark markup is gone.\n$updatePara2 = Find-ParagraphByText \"UpdatePODParameter\"\n$updateIndex = $updatePara2.Range.Start\n\n$emptySeen = 0\n$bookmarkPara = $null\n$passed = $false\nforeach ($p in $d.Paragraphs) {\n    if (-not $passed) {\n        if ($p.Range.Start -eq $updateIndex) { $passed = $true }\n        continue\n    }\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq \"\") {\n        $emptySeen++\n        if ($emptySeen -eq 3) {\n            $bookmarkPara = $p\n            break\n        }\n    } else {\n        break\n    }\n}\n\nif ($null -ne $bookmarkPara) {\n    [void]$bookmarkPara.Range.InsertXML((New-Pkg '<w:p/>'))\n}\n"}
